# Weekly fruit/vegetable price update:
# Insert a new weekly record for "Poroto verde" (Terminal Hortofrutícola Agro
# Chillán) above the existing row 57, shifting the following rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 57 (existing rows 57:63 shift down to 58:64)
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with this week's record
$ws.Range("A57").Value = 7
$ws.Range("B57").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C57").Value = "Ñuble"
$ws.Range("D57").Value = 44585
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = 100112031
$ws.Range("G57").Value = "Poroto verde"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 100
$ws.Range("K57").Value = 25000
$ws.Range("L57").Value = 26000
$ws.Range("M57").Value = 25500
$ws.Range("N57").Value = "$/saco 25 kilos"
$ws.Range("O57").Value = "Región del Maule"
$ws.Range("P57").Value = 1020
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"
